$wb = $excel.ActiveWorkbook

# Sheet ALC Row 10
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()

# Sheet ALC Row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 19241824
$ws.Range("I62").Value = 31264900
$ws.Range("J62").Value = 4901.2
$ws.Range("K62").Value = 31264900
$ws.Range("L62").Value = 4901.2
$ws.Range("M62").Value = -31264276
$ws.Range("N62").Value = -6149.2

# Sheet ALC Row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 19241824
$ws.Range("I65").Value = 31264900
$ws.Range("J65").Value = 4901.2
$ws.Range("K65").Value = 156324500
$ws.Range("L65").Value = 24506
$ws.Range("M65").Value = -156321380
$ws.Range("N65").Value = -30746

# Sheet ALC Row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 372.86365
$ws.Range("I92").Value = 383.9375
$ws.Range("K92").Value = 383.9375
$ws.Range("M92").Value = 864.0625

# Sheet ALC Row 96
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 1245
$ws.Range("I96").Value = 974.5
$ws.Range("J96").Value = 1554.1428
$ws.Range("K96").Value = 2923.5
$ws.Range("L96").Value = 4662.428400000001
$ws.Range("M96").Value = -1550.5
$ws.Range("N96").Value = -7408.428400000001

# Sheet ALC Row 101
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 7742593.5
$ws.Range("I101").Value = 200203.4
$ws.Range("J101").Value = 45454544
$ws.Range("K101").Value = 600610.2
$ws.Range("L101").Value = 136363632
$ws.Range("M101").Value = -598988.2
$ws.Range("N101").Value = -136366876

# Sheet ALC Row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 7467371
$ws.Range("I132").Value = 8932578
$ws.Range("J132").Value = 8133.1816
$ws.Range("K132").Value = 26797734
$ws.Range("L132").Value = 24399.5448
$ws.Range("M132").Value = -26795204
$ws.Range("N132").Value = -29459.5448

# Sheet ARM Row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 72470.07000000001
$ws.Range("I45").Value = 91610.73
$ws.Range("J45").Value = 2287.6667
$ws.Range("K45").Value = 91610.73
$ws.Range("L45").Value = 2287.6667
$ws.Range("M45").Value = -91233.73
$ws.Range("N45").Value = -3041.6667

# Sheet ARM Row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1004.3111
$ws.Range("I61").Value = 982.18604
$ws.Range("J61").Value = 1480
$ws.Range("K61").Value = 982.18604
$ws.Range("L61").Value = 1480
$ws.Range("M61").Value = -770.18604
$ws.Range("N61").Value = -1904

# Sheet ARM Row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1004.3111
$ws.Range("I136").Value = 982.18604
$ws.Range("J136").Value = 1480
$ws.Range("K136").Value = 2946.55812
$ws.Range("L136").Value = 4440
$ws.Range("M136").Value = -396.5581200000001
$ws.Range("N136").Value = -9540

# Sheet BSM Row 25
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 502.5
$ws.Range("I25").Value = 502.5
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 502.5
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -267.5
$ws.Range("N25").ClearContents()

# Sheet CRP Row 12
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 1051.25
$ws.Range("I12").Value = 1051.25
$ws.Range("K12").Value = 1051.25
$ws.Range("M12").Value = -881.25

# Sheet CRP Row 41
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 5500
$ws.Range("I41").Value = 5500
$ws.Range("K41").Value = 5500
$ws.Range("M41").Value = -5072

# Sheet CRP Row 50
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

# Sheet CRP Row 51
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 17913.715
$ws.Range("I51").Value = 7000
$ws.Range("J51").Value = 19732.666
$ws.Range("K51").Value = 7000
$ws.Range("L51").Value = 19732.666
$ws.Range("M51").Value = -6264
$ws.Range("N51").Value = -21204.666

# Sheet CRP Row 59
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 200
$ws.Range("I59").Value = 200
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 200
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = 945
$ws.Range("N59").ClearContents()

# Sheet CRP Row 60
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 20400.889
$ws.Range("I60").Value = 4046.5
$ws.Range("J60").Value = 25073.572
$ws.Range("K60").Value = 4046.5
$ws.Range("L60").Value = 25073.572
$ws.Range("M60").Value = -3535.5
$ws.Range("N60").Value = -26095.572

# Sheet CRP Row 61
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 17913.715
$ws.Range("I61").Value = 7000
$ws.Range("J61").Value = 19732.666
$ws.Range("K61").Value = 7000
$ws.Range("L61").Value = 19732.666
$ws.Range("M61").Value = -6652
$ws.Range("N61").Value = -20428.666

# Sheet CRP Row 68
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 21392.5
$ws.Range("J68").Value = 21392.5
$ws.Range("L68").Value = 21392.5
$ws.Range("N68").Value = -22890.5

# Sheet CRP Row 71
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 21392.5
$ws.Range("J71").Value = 21392.5
$ws.Range("L71").Value = 64177.5
$ws.Range("N71").Value = -71665.5

# Sheet CRP Row 74
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 26314
$ws.Range("J74").Value = 26314
$ws.Range("L74").Value = 26314
$ws.Range("N74").Value = -28062

# Sheet CRP Row 77
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 26314
$ws.Range("J77").Value = 26314
$ws.Range("L77").Value = 78942
$ws.Range("N77").Value = -87678

# Sheet CUL Row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 24409.586
$ws.Range("J12").Value = 32282.322
$ws.Range("L12").Value = 96846.966
$ws.Range("N12").Value = -97192.966

# Sheet GSM Row 22
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 600
$ws.Range("I22").Value = 400.5
$ws.Range("J22").Value = 999
$ws.Range("K22").Value = 400.5
$ws.Range("L22").Value = 999
$ws.Range("M22").Value = 128.5
$ws.Range("N22").Value = -2057

# Sheet LTW Row 5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 3500
$ws.Range("I5").Value = 2000
$ws.Range("J5").Value = 3875
$ws.Range("K5").Value = 2000
$ws.Range("L5").Value = 3875
$ws.Range("M5").Value = -1887
$ws.Range("N5").Value = -4101

# Sheet LTW Row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2452678.2
$ws.Range("I46").Value = 5952978.5
$ws.Range("J46").Value = 2468
$ws.Range("K46").Value = 5952978.5
$ws.Range("L46").Value = 2468
$ws.Range("M46").Value = -5952790.5
$ws.Range("N46").Value = -2844

# Sheet WVR Row 22
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

# Sheet WVR Row 26
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 3000
$ws.Range("I26").Value = 3000
$ws.Range("K26").Value = 3000
$ws.Range("M26").Value = -2707

# Sheet WVR Row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 17400.312
$ws.Range("I122").Value = 27604.79
$ws.Range("K122").Value = 82814.37
$ws.Range("M122").Value = -80364.37
